# Fruta / hortaliza, semanal
# Insert a new weekly data row before the current row 10 (Fecha=2022-06-08 / serial 44720),
# pushing the existing rows 10-21 down to rows 11-22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 10..21 down to 11..22, duplicating row 10's formatting into the new row.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly record.
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "Macroferia Regional de Talca"
$ws.Range("C10").Value = "Maule"
$ws.Range("D10").Value = 44720
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 100112040
$ws.Range("G10").Value = "Cilantro"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 150
$ws.Range("K10").Value = 9000
$ws.Range("L10").Value = 9000
$ws.Range("M10").Value = 9000
$ws.Range("N10").Value = "$/caja 36 atados"
$ws.Range("O10").Value = "Región Metropolitana"
$ws.Range("P10").Value = 250
$ws.Range("Q10").Value = 36
$ws.Range("R10").Value = "Hortaliza"

# Make sure the date cell keeps the workbook's date number format (style s="2").
$ws.Range("D10").NumberFormat = $ws.Range("D11").NumberFormat
